# Update "想去人数" (interest count) values in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 3132
$ws1.Range("F7").Value = 3904
$ws1.Range("F8").Value = 480

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 3132
$ws4.Range("F8").Value = 3904
$ws4.Range("F9").Value = 480
